# Edits "Architecture de code" document:
#  - wraps code/technical terms in w:proofErr spellStart/spellEnd markers
#    (mirrors Word's automatic spell-checker markup around runs that were
#    retouched)
#  - rewrites the "Interactions" bullet wording
#  - appends a new "State" bullet with three sub-bullets (Time.cs, Life.cs,
#    Inventary.cs)
#
# Note: Range.InsertXML() on this runtime always splices the new content in
# at the *end* of the paragraph being edited rather than strictly in place,
# so every replacement below targets a whole paragraph (Start .. End-1, i.e.
# everything except the trailing pilcrow) and reproduces any trailing runs
# that must be kept unchanged.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$nbsp = [char]0x00A0
$apos = [char]0x2019
$eacute = [char]0x00E9

function Find-ParagraphStartingWith([string]$prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($prefix)) {
            return $d.Paragraphs($i)
        }
    }
    throw "No paragraph starting with: $prefix"
}

function Replace-ParagraphBody([string]$prefix, [string]$innerXml) {
    # Replace the full body (excluding the trailing pilcrow) of the
    # paragraph that starts with $prefix, with $innerXml.
    $p = Find-ParagraphStartingWith $prefix
    $s = $p.Range.Start
    $e = $p.Range.End - 1
    $r = $d.Range($s, $e)
    $xml = $pkgHeader + '<w:p>' + $innerXml + '</w:p>' + $pkgFooter
    $ignored = $r.InsertXML($xml)
}

# 1. "Characters (nous mettrons tout ce qui concerne les personnages)"
Replace-ParagraphBody "Characters (nous mettrons" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Characters</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (nous mettrons tout ce qui concerne les personnages)</w:t></w:r>'
)

# 2. "Characters.cs" -> wrap in proofErr
Replace-ParagraphBody "Characters.cs" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Characters.cs</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# 3. "Interactions (...)" paragraph rewrite (keeps the existing bookmark)
Replace-ParagraphBody "Interactions" (
    '<w:r><w:t>Interactions</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (Tout ce qui concerne </w:t></w:r>' +
    '<w:r><w:t>les interactions</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve"> du personnage avec le d' + $eacute + 'cor et les objets</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>'
)

# 4. "Collisions.cs" -> wrap in proofErr
Replace-ParagraphBody "Collisions.cs" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Collisions.cs</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# 5. "Recuperate.cs" -> wrap in proofErr
Replace-ParagraphBody "Recuperate.cs" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Recuperate.cs</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# 6. "Moves.cs" -> wrap in proofErr
Replace-ParagraphBody "Moves.cs" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Moves.cs</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# 7. "Card (la carte du jeu)"
Replace-ParagraphBody "Card (la carte du jeu)" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Card</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> (la carte du jeu)</w:t></w:r>'
)

# 8. "Card.cs ???" -> wrap "Card.cs" in proofErr, keep the trailing " ???" run
Replace-ParagraphBody "Card.cs" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Card.cs</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">' + $nbsp + '???</w:t></w:r>'
)

# 9. "CardDesign.cs" -> wrap in proofErr
Replace-ParagraphBody "CardDesign.cs" (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>CardDesign.cs</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)

# 10. Append new "State" bullet with Time.cs / Life.cs / Inventary.cs sub-bullets
$lastP = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastP.Range.End, $lastP.Range.End)

$newParas =
    '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t>State (tout ce qui concerne l' + $apos + $eacute + 'tat du personnage)</w:t></w:r></w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Time.cs</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Life.cs</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Inventary.cs</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

$ignored = $insertPoint.InsertXML($pkgHeader + $newParas + $pkgFooter)
